$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update probability values in row 8
$ws.Range("K8").Value = 0.5
$ws.Range("L8").Value = 0.5
$ws.Range("N8").Value = 0

# Update the active cell selection on the sheet
$ws.Activate()
$ws.Range("Q8").Select()
